$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the crypto price/volume refresh described by the commit.
# A handful of Price (column D) values are plain decimal numbers
# (e.g. '0.9987', '1.000') that Excel's normal type-inference would
# silently convert to numeric cells; the source data stores them as
# literal text, so those specific cells are pinned to text format
# ('@') before the value is written. Values that are already
# unambiguous as text (names, URLs, multi-dot prices, padded percents)
# are written directly, leaving their formatting untouched.

$ws.Range('D2').Value = '24.121.70'
$ws.Range('E2').Value = '  -3.04%  '
$ws.Range('D3').Value = '1.640.12'
$ws.Range('E3').Value = '  -2.92%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9987'
$ws.Range('E4').Value = '  -0.69%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '307.39'
$ws.Range('E5').Value = '  -2.35%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$ws.Range('E6').Value = '  -0.46%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.3883'
$ws.Range('E7').Value = '  -1.47%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3852'
$ws.Range('E8').Value = '  -3.36%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.9994'
$ws.Range('E9').Value = '  -0.49%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '49.44'
$ws.Range('E10').Value = '  -5.77%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.347'
$ws.Range('E11').Value = '  -6.52%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.08743'
$ws.Range('E12').Value = '  +0.35%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '23.65'
$ws.Range('E13').Value = '  -6.91%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '7.082'
$ws.Range('E14').Value = '  -3.84%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.00001286'
$ws.Range('E15').Value = '  -3.25%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '7.450'
$ws.Range('E16').Value = '  -5.19%  '
$ws.Range('D17').Value = '1.632.54'
$ws.Range('E17').Value = '  -2.03%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '95.28'
$ws.Range('E18').Value = '  +0.67%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06883'
$ws.Range('E19').Value = '  -4.17%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '20.80'
$ws.Range('E20').Value = '  +2.09%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.873'
$ws.Range('E21').Value = '  -4.10%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '1.000'
$ws.Range('E22').Value = '  -0.43%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '13.51'
$ws.Range('E23').Value = '  -4.56%  '
$ws.Range('D24').Value = '24.107.87'
$ws.Range('E24').Value = '  -3.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.328'
$ws.Range('E25').Value = '  -3.22%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.695'
$ws.Range('E26').Value = '  -5.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '22.31'
$ws.Range('E27').Value = '  -2.95%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '157.14'
$ws.Range('E28').Value = '  -3.26%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '8.679'
$ws.Range('E29').Value = '  +8.26%  '
$ws.Range('B30').Value = 'HuobiToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '5.365'
$ws.Range('E30').Value = '  -11.19%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '139.99'
$ws.Range('E31').Value = '  -5.51%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.414'
$ws.Range('E32').Value = '  -9.20%  '
$ws.Range('D33').Value = '1.815.15'
$ws.Range('E33').Value = '  -7.67%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '6.842'
$ws.Range('E34').Value = '  -2.16%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.08004'
$ws.Range('E35').Value = '  -6.12%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.02877'
$ws.Range('E36').Value = '  -7.49%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.2665'
$ws.Range('E37').Value = '  -6.84%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.9453'
$ws.Range('E38').Value = '  -8.33%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.09177'
$ws.Range('E39').Value = '  -5.25%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.445'
$ws.Range('E40').Value = '  -1.62%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '9.855'
$ws.Range('E41').Value = '  -8.26%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.7521'
$ws.Range('E42').Value = '  -6.69%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '13.03'
$ws.Range('E43').Value = '  -5.81%  '
$ws.Range('E44').Value = '  -6.17%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.6868'
$ws.Range('E45').Value = '  -5.35%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.453'
$ws.Range('E46').Value = '  -6.47%  '
$ws.Range('E47').Value = '  -3.25%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.9991'
$ws.Range('E48').Value = '  -0.66%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.08359'
$ws.Range('E49').Value = '  -6.36%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.256'
$ws.Range('E50').Value = '  -9.26%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '132.17'
$ws.Range('E51').Value = '  -4.84%  '
